{"js": "// The document starts with two inline screenshots (paragraphs 1 and 2),\n// followed by the \"History cost\" and \"Nonlinear Cost function\" paragraphs.\n// This edit removes the second screenshot - the one sitting directly above\n// \"History cost\" - while leaving the first screenshot and the text intact.\n// The paragraph that held the picture is left in place, now empty.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// Load each paragraph's inline pictures so we can find the right one.\nconst pictureCollections = paragraphs.items.map((p) => {\n  const pics = p.inlinePictures;\n  pics.load(\"items\");\n  return pics;\n});\nawait context.sync();\n\n// Collect every paragraph index that currently has an inline picture.\nconst pictureParagraphIndexes = [];\nfor (let i = 0; i < pictureCollections.length; i++) {\n  if (pictureCollections[i].items.length > 0) {\n    pictureParagraphIndexes.push(i);\n  }\n}\n\n// The picture to remove is the second picture in the document (the one\n// immediately preceding the \"History cost\" paragraph) - i.e. the last of\n// the two leading picture paragraphs.\nif (pictureParagraphIndexes.length >= 2) {\n  const targetIndex = pictureParagraphIndexes[1];\n  const targetPictures = pictureCollections[targetIndex];\n  targetPictures.items[0].delete();\n  await context.sync();\n}\n", "ps1": "# The document starts with two inline screenshots (the first two paragraphs),\n# followed by the \"History cost\" and \"Nonlinear Cost function\" paragraphs.\n# This edit removes the second screenshot - the one sitting directly above\n# \"History cost\" - while leaving the first screenshot and the text intact.\n# The paragraph that held the picture is left in place, now empty.\n\n$d = $word.ActiveDocument\n\n# wdInlineShapePicture = 3\n$wdInlineShapePicture = 3\n\n$pictureShapes = @()\nfor ($i = 1; $i -le $d.InlineShapes.Count; $i++) {\n    $shape = $d.InlineShapes.Item($i)\n    if ($shape.Type -eq $wdInlineShapePicture) {\n        $pictureShapes += $shape\n    }\n}\n\n# The picture to remove is the second picture in the document (the one\n# immediately preceding the \"History cost\" paragraph).\nif ($pictureShapes.Count -ge 2) {\n    $pictureShapes[1].Delete()\n}\n"}
